$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "general" (tab 1): add two missing results, then make this the
# selected/active sheet with G8 selected.
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Range("G4").Value = 3.4
$wsGeneral.Range("G7").Value = 3.3

# ---------------------------------------------------------------------------
# Sheet "pretraining" (tab 2): selection only moves, no data change.
# ---------------------------------------------------------------------------
$wsPretraining = $wb.Worksheets.Item("pretraining")

# ---------------------------------------------------------------------------
# Sheet "window_size" (tab 3): add missing results.
# ---------------------------------------------------------------------------
$wsWindowSize = $wb.Worksheets.Item("window_size")
$wsWindowSize.Range("F4").Value = 7.8
$wsWindowSize.Range("E5").Value = 3.3
$wsWindowSize.Range("F5").Value = 7.6
$wsWindowSize.Range("F8").Value = 7.7

# ---------------------------------------------------------------------------
# Sheet "scf_size" (tab 4): add missing results.
# ---------------------------------------------------------------------------
$wsScfSize = $wb.Worksheets.Item("scf_size")
$wsScfSize.Range("G5").Value = 3.4
$wsScfSize.Range("H6").Value = 7.8
$wsScfSize.Range("G9").Value = 3.4

# ---------------------------------------------------------------------------
# Sheet "w2v_size" (tab 5): add many missing results, fix one existing one.
# ---------------------------------------------------------------------------
$wsW2vSize = $wb.Worksheets.Item("w2v_size")
$wsW2vSize.Range("F4").Value = 3.3
$wsW2vSize.Range("F5").Value = 3.3
$wsW2vSize.Range("F6").Value = 3.4
$wsW2vSize.Range("F8").Value = 3.5
$wsW2vSize.Range("G8").Value = 7.8
$wsW2vSize.Range("F9").Value = 3.4
$wsW2vSize.Range("F10").Value = 3.5
$wsW2vSize.Range("F11").Value = 3.3
$wsW2vSize.Range("F12").Value = 3.4
$wsW2vSize.Range("G12").Value = 7.8
$wsW2vSize.Range("F13").Value = 3.3
$wsW2vSize.Range("F14").Value = 3.5
$wsW2vSize.Range("F15").Value = 3.5
$wsW2vSize.Range("F16").Value = 3.4
$wsW2vSize.Range("G16").Value = 8.2

# ---------------------------------------------------------------------------
# Sheet "w2v_proj" (tab 6): add missing results.
# ---------------------------------------------------------------------------
$wsW2vProj = $wb.Worksheets.Item("w2v_proj")
$wsW2vProj.Range("E4").Value = 3.3
$wsW2vProj.Range("F4").Value = 7.4
$wsW2vProj.Range("E5").Value = 3.4

# ---------------------------------------------------------------------------
# Update each sheet's remembered selection (activeCell) to match the final
# state, then finish with "general" activated (and its selection on G8),
# which becomes the workbook's active tab.
# ---------------------------------------------------------------------------
$wsPretraining.Range("E7").Select() | Out-Null
$wsWindowSize.Range("E6").Select() | Out-Null
$wsScfSize.Range("G10").Select() | Out-Null
$wsW2vSize.Range("F14").Select() | Out-Null
$wsW2vProj.Range("E6").Select() | Out-Null

$wsGeneral.Activate() | Out-Null
$wsGeneral.Range("G8").Select() | Out-Null

Write-Host "edits applied"
